# Apply the changes described by the diff to the active workbook.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Лист1")

# --- Title cells on row 1 (drop the period right after "6.4.1.2") ---
# A1 (Kyrgyz title) text is unchanged, only B1/C1 text actually changes.
$ws.Range("B1").Value = "6.4.1.2 Потери воды при транспортировке"
$ws.Range("C1").Value = "6.4.1.2 Percentage of water loss during transportation"

# --- Updated data values ---
$ws.Range("P5").Value = 2388
$ws.Range("P10").Value = 335.3
$ws.Range("P16").Value = 27.3
$ws.Range("P21").Value = 24.3

# --- Selection moves from R9 to S3 ---
$ws.Range("S3").Select()
